# Daily attendance processing - reorders the "Recorded By" (column G)
# comma-separated contributor list on the active sheet: each multi-value
# entry is reversed (last contributor becomes first, etc.).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ur = $ws.UsedRange
$lastRow = $ur.Rows.Count()

$changedCount = 0
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()
    if ($val -eq $null) { continue }

    $parts = $val.Split(",")
    if ($parts.Count -le 1) { continue }

    $n = $parts.Count
    $rev = @()
    for ($i = $n - 1; $i -ge 0; $i--) {
        $rev += $parts[$i].Trim()
    }
    $newVal = [string]::Join(", ", $rev)

    $cell.Value = $newVal
    $changedCount++
}

Write-Host "Reordered Recorded By values in $changedCount rows"
